{"js": "// Insert a new numbered-list paragraph right after the paragraph that ends\n// with \"...which user added the peep.\" The new paragraph reads:\n// \"Tried everything. Didn't work. Now wondering what's wrong with this code of mine. \"\n// and keeps the same list formatting (ListParagraph style, same numbering).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst marker = \"going to add a tag for each user\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the target paragraph.\");\n}\n\nconst newText =\n  \"Tried everything. Didn\\u2019t work. Now wondering what\\u2019s wrong with this code of mine. \";\n\n// insertParagraph(\"After\") creates a sibling paragraph that inherits the\n// source paragraph's formatting (style, numbering, run properties), which is\n// exactly the ListParagraph / numId=1 bullet formatting we need here.\ntarget.insertParagraph(newText, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert a new numbered-list paragraph right after the paragraph that ends\n# with \"...which user added the peep.\" The new paragraph reads:\n# \"Tried everything. Didn't work. Now wondering what's wrong with this code of mine. \"\n# and keeps the same list formatting (ListParagraph style, same numbering),\n# because InsertParagraphAfter() on the existing list paragraph's Range\n# inherits that paragraph's formatting for the freshly created paragraph.\n\n$rsq = [char]0x2019\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"*going to add a tag for each user*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the target paragraph.\"\n}\n\n$newText = \"Tried everything. Didn\" + $rsq + \"t work. Now wondering what\" + $rsq + \"s wrong with this code of mine. \"\n\n$target = $d.Paragraphs($targetIndex)\n$target.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs($targetIndex + 1)\n$newPara.Range.Text = $newText\n"}
